$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 230
$ws.Range("I8").Value = 230
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 690
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -551
$ws.Range("N8").ClearContents()
$ws.Range("H17").Value = 3995
$ws.Range("J17").Value = 3995
$ws.Range("L17").Value = 11985
$ws.Range("N17").Value = -12321
$ws.Range("H31").Value = 249.5
$ws.Range("I31").Value = 249.5
$ws.Range("K31").Value = 748.5
$ws.Range("M31").Value = -518.5
$ws.Range("H33").Value = 174.2
$ws.Range("I33").Value = 171.71428
$ws.Range("K33").Value = 171.71428
$ws.Range("M33").Value = 57.28572
$ws.Range("H38").Value = 6399.8
$ws.Range("I38").Value = 999.5
$ws.Range("J38").Value = 10000
$ws.Range("K38").Value = 2998.5
$ws.Range("L38").Value = 30000
$ws.Range("M38").Value = -2626.5
$ws.Range("N38").Value = -30744
$ws.Range("H80").Value = 3088.111
$ws.Range("I80").Value = 1649.5
$ws.Range("J80").Value = 3499.1428
$ws.Range("K80").Value = 4948.5
$ws.Range("L80").Value = 10497.4284
$ws.Range("M80").Value = -3950.5
$ws.Range("N80").Value = -12493.4284
$ws.Range("H83").Value = 3088.111
$ws.Range("I83").Value = 1649.5
$ws.Range("J83").Value = 3499.1428
$ws.Range("K83").Value = 14845.5
$ws.Range("L83").Value = 31492.2852
$ws.Range("M83").Value = -9853.5
$ws.Range("N83").Value = -41476.2852
$ws.Range("I92").Value = 985.1177
$ws.Range("J92").Value = 1033.1666
$ws.Range("K92").Value = 985.1177
$ws.Range("L92").Value = 1033.1666
$ws.Range("M92").Value = 262.8823
$ws.Range("N92").Value = -3529.1666
$ws.Range("H96").Value = 1369.909
$ws.Range("I96").Value = 341.1111
$ws.Range("K96").Value = 1023.3333
$ws.Range("M96").Value = 349.6667
$ws.Range("I101").Value = 485
$ws.Range("J101").Value = 650
$ws.Range("K101").Value = 1455
$ws.Range("L101").Value = 1950
$ws.Range("M101").Value = 167
$ws.Range("N101").Value = -5194
$ws.Range("H113").Value = 4275.25
$ws.Range("I113").Value = 3675.5
$ws.Range("J113").Value = 4875
$ws.Range("K113").Value = 3675.5
$ws.Range("L113").Value = 4875
$ws.Range("M113").Value = -421.5
$ws.Range("N113").Value = -11383
$ws.Range("H132").Value = 1507.5161
$ws.Range("I132").Value = 1507.5161
$ws.Range("K132").Value = 4522.5483
$ws.Range("M132").Value = -1992.5483

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 9531.714
$ws.Range("I74").Value = 9531.714
$ws.Range("K74").Value = 9531.714
$ws.Range("M74").Value = -8657.714
$ws.Range("H77").Value = 9531.714
$ws.Range("I77").Value = 9531.714
$ws.Range("K77").Value = 47658.57
$ws.Range("M77").Value = -43290.57
$ws.Range("H110").Value = 3058.7144
$ws.Range("I110").Value = 1367.4166
$ws.Range("J110").Value = 13206.5
$ws.Range("K110").Value = 1367.4166
$ws.Range("L110").Value = 13206.5
$ws.Range("M110").Value = 677.5834
$ws.Range("N110").Value = -17296.5
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 9533.111000000001
$ws.Range("I20").Value = 6966.5
$ws.Range("K20").Value = 6966.5
$ws.Range("M20").Value = -6719.5
$ws.Range("H44").Value = 60000
$ws.Range("J44").Value = 60000
$ws.Range("L44").Value = 60000
$ws.Range("N44").Value = -60994
$ws.Range("H107").Value = 2043.5555
$ws.Range("I107").Value = 1898.8334
$ws.Range("J107").Value = 2333
$ws.Range("K107").Value = 1898.8334
$ws.Range("L107").Value = 2333
$ws.Range("M107").Value = 21.16660000000002
$ws.Range("N107").Value = -6173

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H82").Value = 10000
$ws.Range("J82").Value = 10000
$ws.Range("L82").Value = 10000
$ws.Range("N82").Value = -10722
$ws.Range("H85").Value = 10000
$ws.Range("J85").Value = 10000
$ws.Range("L85").Value = 10000
$ws.Range("N85").Value = -12496
$ws.Range("H122").Value = 4101.0713
$ws.Range("I122").Value = 3326.75
$ws.Range("K122").Value = 9980.25
$ws.Range("M122").Value = -7530.25
$ws.Range("H141").Value = 227854.28
$ws.Range("J141").Value = 227854.28
$ws.Range("L141").Value = 227854.28
$ws.Range("N141").Value = -238214.28

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1013943.8
$ws.Range("I4").Value = 2534709.5
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 7604128.5
$ws.Range("L4").Value = 300
$ws.Range("M4").Value = -7604016.5
$ws.Range("N4").Value = -524
$ws.Range("H5").Value = 662.1667
$ws.Range("I5").Value = 458
$ws.Range("J5").Value = 866.3333
$ws.Range("K5").Value = 1374
$ws.Range("L5").Value = 2598.9999
$ws.Range("M5").Value = -1262
$ws.Range("N5").Value = -2822.9999
$ws.Range("H23").Value = 50
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H34").Value = 7547
$ws.Range("J34").Value = 7547
$ws.Range("L34").Value = 22641
$ws.Range("N34").Value = -22809
$ws.Range("H109").Value = 4229.8696
$ws.Range("I109").Value = 2822
$ws.Range("J109").Value = 4526.263
$ws.Range("K109").Value = 8466
$ws.Range("L109").Value = 13578.789
$ws.Range("M109").Value = -7426
$ws.Range("N109").Value = -15658.789
$ws.Range("H129").Value = 1094.25
$ws.Range("I129").Value = 1094.25
$ws.Range("K129").Value = 3282.75
$ws.Range("M129").Value = 1717.25
$ws.Range("H131").Value = 3138.3142
$ws.Range("J131").Value = 3161.0322
$ws.Range("L131").Value = 9483.096600000001
$ws.Range("N131").Value = -19563.0966
$ws.Range("H135").Value = 662.1667
$ws.Range("I135").Value = 458
$ws.Range("J135").Value = 866.3333
$ws.Range("K135").Value = 4122
$ws.Range("L135").Value = 7796.9997
$ws.Range("M135").Value = -1587
$ws.Range("N135").Value = -12866.9997

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H132").Value = 4038.6365
$ws.Range("I132").Value = 3571.5
$ws.Range("K132").Value = 10714.5
$ws.Range("M132").Value = -8184.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3494.3333
$ws.Range("I7").Value = 3494.3333
$ws.Range("K7").Value = 3494.3333
$ws.Range("M7").Value = -3382.3333
$ws.Range("H61").Value = 1301.3334
$ws.Range("I61").Value = 952
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 952
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -750
$ws.Range("N61").Value = -2404
$ws.Range("H93").Value = 2089.889
$ws.Range("I93").Value = 1929.3334
$ws.Range("J93").Value = 2250.4443
$ws.Range("K93").Value = 1929.3334
$ws.Range("L93").Value = 2250.4443
$ws.Range("M93").Value = -681.3334
$ws.Range("N93").Value = -4746.4443
$ws.Range("H113").Value = 1301.3334
$ws.Range("I113").Value = 952
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 952
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 1218
$ws.Range("N113").Value = -6340
$ws.Range("H126").Value = 3494.3333
$ws.Range("I126").Value = 3494.3333
$ws.Range("K126").Value = 10482.9999
$ws.Range("M126").Value = -8012.999899999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8291.857
$ws.Range("J62").Value = 8940.5
$ws.Range("L62").Value = 8940.5
$ws.Range("N62").Value = -10188.5
$ws.Range("H65").Value = 8291.857
$ws.Range("J65").Value = 8940.5
$ws.Range("L65").Value = 44702.5
$ws.Range("N65").Value = -50942.5
$ws.Range("H81").Value = 2542.2727
$ws.Range("I81").Value = 1896.5
$ws.Range("J81").Value = 9000
$ws.Range("K81").Value = 3793
$ws.Range("L81").Value = 18000
$ws.Range("M81").Value = -2732
$ws.Range("N81").Value = -20122
$ws.Range("H84").Value = 2542.2727
$ws.Range("I84").Value = 1896.5
$ws.Range("J84").Value = 9000
$ws.Range("K84").Value = 18965
$ws.Range("L84").Value = 90000
$ws.Range("M84").Value = -13661
$ws.Range("N84").Value = -100608
$ws.Range("I100").Value = 749.5714
$ws.Range("J100").Value = 854
$ws.Range("K100").Value = 1499.1428
$ws.Range("L100").Value = 1708
$ws.Range("M100").Value = -958.1428000000001
$ws.Range("N100").Value = -2790
$ws.Range("H107").Value = 731.9167
$ws.Range("I107").Value = 597.8570999999999
$ws.Range("J107").Value = 919.6
$ws.Range("K107").Value = 1793.5713
$ws.Range("L107").Value = 2758.8
$ws.Range("M107").Value = 126.4287000000002
$ws.Range("N107").Value = -6598.8
$ws.Range("H132").Value = 6542.143
$ws.Range("I132").Value = 2699.5
$ws.Range("K132").Value = 8098.5
$ws.Range("M132").Value = -5568.5
$ws.Range("H136").Value = 13750
$ws.Range("I136").Value = 13333.333
$ws.Range("J136").Value = 15000
$ws.Range("K136").Value = 39999.999
$ws.Range("L136").Value = 45000
$ws.Range("M136").Value = -37449.999
$ws.Range("N136").Value = -50100
